$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 996
$ws1.Range("F9").Value = 1469
$ws1.Range("F11").Value = 1369
$ws1.Range("F12").Value = 3024
$ws1.Range("F13").Value = 476
$ws1.Range("F14").Value = 1658
$ws1.Range("F18").Value = 1406
$ws1.Range("F21").Value = 1137
$ws1.Range("F22").Value = 6
$ws1.Range("F23").Value = 407
$ws1.Range("F25").Value = 3530
$ws1.Range("F26").Value = 700

# Sheet "全部类型" (sheet4) - column F "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 996
$ws4.Range("F19").Value = 1469
$ws4.Range("F21").Value = 1369
$ws4.Range("F22").Value = 3024
$ws4.Range("F23").Value = 476
$ws4.Range("F24").Value = 1658
$ws4.Range("F28").Value = 1406
$ws4.Range("F33").Value = 1137
$ws4.Range("F34").Value = 6
$ws4.Range("F35").Value = 407
$ws4.Range("F37").Value = 3530
$ws4.Range("F38").Value = 700

$wb.Save()
